$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Team Name (col B) and Users (col C) for rows 2..21 (teams 1..20),
# and reset Score (col J) to 100 for each of those rows.
for ($i = 1; $i -le 20; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = "Team $i"

    $u1 = 4 * ($i - 1) + 1
    $u2 = $u1 + 1
    $u3 = $u1 + 2
    $u4 = $u1 + 3
    $ws.Cells.Item($row, 3).Value = "User$u1, User$u2, User$u3, User$u4"

    $ws.Cells.Item($row, 10).Value = 100
}

# Row 16 (team 15 / formerly "BABLU") no longer carries the "5, 2, 7" Powerups value.
$ws.Cells.Item(16, 8).Value = ""

# Remove the trailing teams (rows 22..28, serial numbers 21..27) entirely.
$ws.Range("A22:J28").EntireRow.Delete()
